# Prefix the header row (row 1) labels with "Ano" / "Intervalo" so that Power BI
# can automatically promote the first row to a header.

$wb = $excel.ActiveWorkbook

# Sheets 1, 2, 3, 5: simple "Ano <year>" labels for B1:E1
$anoSheets = @(1, 2, 3, 5)
foreach ($idx in $anoSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("B1").Value = "Ano 2015"
    $ws.Range("C1").Value = "Ano 2030"
    $ws.Range("D1").Value = "Ano 2040"
    $ws.Range("E1").Value = "Ano 2050"
}

# Sheet 4: "Intervalo <range>" labels for B1:E1
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B1").Value = "Intervalo 2015"
$ws4.Range("C1").Value = "Intervalo 2015-2030"
$ws4.Range("D1").Value = "Intervalo 2031-2040"
$ws4.Range("E1").Value = "Intervalo 2041-2050"

# Sheet 6: only has B1 -> "Ano 2015"
$ws6 = $wb.Worksheets.Item(6)
$ws6.Range("B1").Value = "Ano 2015"
